$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# Row 8: add Description for logFileTemplatePath
$ws.Range("C8").Value = "ไฟล์ Template ที่มี Header ของ Log"

# Row 9: new entry blankTemplatePath
$ws.Range("A9").Value = "blankTemplatePath"
$ws.Range("B9").Value = $ws.Range("B8").Value
$ws.Range("B9").Style = $ws.Range("B8").Style
$ws.Range("C9").Value = "ไฟล์ Template เปล่าๆ ไว้ใส่ข้อมูลที่ Log แล้ว"

# Row 10: new entry logPath (header-style row like row 1)
$ws.Range("A10").Value = "logPath"
$ws.Range("B10").Value = "D:\Mean\UIpath Workspace\Output\log\"
$ws.Range("C10").Value = "ไฟล์เก็บ Log"

$ws.Range("B12").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
